$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price (D) and date-like (I) columns stay as literal text,
# not auto-converted to currency numbers / date serials by Excel's smart
# value parsing.
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("I2:I7").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "LOT"
$ws.Range("B2").Value = "1:55 pm"
$ws.Range("C2").Value = "8:20 pm"
$ws.Range("D2").Value = "$590"
$ws.Range("E2").Value = "1 stop (WAW)"
$ws.Range("F2").Value = "13h 25m"
$ws.Range("I2").Value = "2025-04-16"
$ws.Range("K2").Value = 4

# Row 3
$ws.Range("A3").Value = "LOT"
$ws.Range("B3").Value = "1:55 pm"
$ws.Range("C3").Value = "8:20 pm"
$ws.Range("D3").Value = "$558"
$ws.Range("E3").Value = "1 stop (WAW)"
$ws.Range("F3").Value = "13h 25m"
$ws.Range("G3").Value = "IST"
$ws.Range("I3").Value = "2025-04-16"
$ws.Range("K3").Value = 4

# Row 4
$ws.Range("B4").Value = "11:25 am"
$ws.Range("D4").Value = "$568"
$ws.Range("E4").Value = "2 stops (KRK, WAW)"
$ws.Range("F4").Value = "15h 55m"
$ws.Range("I4").Value = "2025-04-16"
$ws.Range("K4").Value = 4

# Row 5
$ws.Range("A5").Value = "LOT"
$ws.Range("B5").Value = "1:55 pm"
$ws.Range("C5").Value = "8:20 pm"
$ws.Range("D5").Value = "$590"
$ws.Range("E5").Value = "1 stop (WAW)"
$ws.Range("F5").Value = "13h 25m"
$ws.Range("I5").Value = "2025-04-17"
$ws.Range("K5").Value = 4

# Row 6
$ws.Range("A6").Value = "LOT"
$ws.Range("B6").Value = "1:55 pm"
$ws.Range("C6").Value = "8:20 pm"
$ws.Range("D6").Value = "$558"
$ws.Range("E6").Value = "1 stop (WAW)"
$ws.Range("F6").Value = "13h 25m"
$ws.Range("I6").Value = "2025-04-17"
$ws.Range("K6").Value = 4

# Row 7 (new)
$ws.Range("A7").Value = "LOT"
$ws.Range("B7").Value = "11:25 am"
$ws.Range("C7").Value = "8:20 pm"
$ws.Range("D7").Value = "$568"
$ws.Range("E7").Value = "2 stops (KRK, WAW)"
$ws.Range("F7").Value = "15h 55m"
$ws.Range("G7").Value = "IST"
$ws.Range("H7").Value = "YYZ"
$ws.Range("I7").Value = "2025-04-17"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = "{}"
